$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7 content corrections:
#  A7: "The_Big_League" -> "The_Big_leauge"
#  D7: password mismatch message gains a trailing period
$ws.Range("A7").Value = "The_Big_leauge"
$ws.Range("D7").Value = "password_mismatch:The two password fields didn" + [char]0x2019 + "t match."

# Update the active selection to D7 (was C20)
$ws.Range("D7").Select()
